$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.216.63"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.53%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.84"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +0.17%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.24%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.76"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("E6").Value = "  +0.25%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5205"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +0.59%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3770"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +0.30%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07273"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +0.82%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.18"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +0.30%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9008"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +0.37%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08191"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +6.96%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.64"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +2.41%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.900.99"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +1.04%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.287"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +1.11%  "

$ws.Range("E16").Value = "  +0.25%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008594"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.95%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.54"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("E19").Value = "  +0.22%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.239.65"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +0.41%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.090"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.61%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +1.16%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.400"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.23%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.310"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +1.18%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.38"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +0.73%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.23"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +0.97%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.743"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.60%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.30"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +0.80%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.822"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +0.70%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.903"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -1.38%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09234"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +0.41%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05041"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -0.05%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7961"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +2.85%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.228"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -0.68%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.433"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +4.68%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.963"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -0.72%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.588"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.33%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5663"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +1.08%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01988"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -0.06%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.073"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.07%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.961"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -0.20%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.563"
$ws.Range("D42").Style = $style

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.14"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -3.34%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1514"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -0.01%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4870"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +1.07%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.07"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -0.61%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.620"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +1.63%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.15"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +2.09%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.48"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -0.63%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05937"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.30%  "
